$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2000
$ws.Range("J43").Value = 2000
$ws.Range("L43").Value = 2000
$ws.Range("N43").Value = -2138
$ws.Range("H70").Value = 2715.6875
$ws.Range("I70").Value = 3000.3333
$ws.Range("J70").Value = 2544.9
$ws.Range("K70").Value = 9000.999899999999
$ws.Range("L70").Value = 7634.700000000001
$ws.Range("M70").Value = -8730.999899999999
$ws.Range("N70").Value = -8174.700000000001
$ws.Range("H73").Value = 2715.6875
$ws.Range("I73").Value = 3000.3333
$ws.Range("J73").Value = 2544.9
$ws.Range("K73").Value = 9000.999899999999
$ws.Range("L73").Value = 7634.700000000001
$ws.Range("M73").Value = -8064.999899999999
$ws.Range("N73").Value = -9506.700000000001
$ws.Range("H106").Value = 83336230
$ws.Range("J106").Value = 166669170
$ws.Range("L106").Value = 166669170
$ws.Range("N106").Value = -166670432
$ws.Range("H107").Value = 8621259
$ws.Range("I107").Value = 12500285
$ws.Range("J107").Value = 1199.8889
$ws.Range("K107").Value = 12500285
$ws.Range("L107").Value = 1199.8889
$ws.Range("M107").Value = -12498365
$ws.Range("N107").Value = -5039.8889
$ws.Range("H116").Value = 6220.64
$ws.Range("I116").Value = 7856.875
$ws.Range("J116").Value = 3311.7778
$ws.Range("K116").Value = 7856.875
$ws.Range("L116").Value = 3311.7778
$ws.Range("M116").Value = -4414.875
$ws.Range("N116").Value = -10195.7778
$ws.Range("H129").Value = 950.95123
$ws.Range("I129").Value = 649.2222
$ws.Range("J129").Value = 988.1507
$ws.Range("K129").Value = 1947.6666
$ws.Range("L129").Value = 2964.4521
$ws.Range("M129").Value = 3052.3334
$ws.Range("N129").Value = -12964.4521
$ws.Range("H132").Value = 1102.7222
$ws.Range("I132").Value = 834.15625
$ws.Range("J132").Value = 3251.25
$ws.Range("K132").Value = 2502.46875
$ws.Range("L132").Value = 9753.75
$ws.Range("M132").Value = 27.53125
$ws.Range("N132").Value = -14813.75
$ws.Range("H135").Value = 1308.0182
$ws.Range("I135").Value = 1154.5
$ws.Range("K135").Value = 10390.5
$ws.Range("M135").Value = -7855.5
$ws.Range("H138").Value = 2442.7537
$ws.Range("I138").Value = 954.73914
$ws.Range("J138").Value = 5418.7827
$ws.Range("K138").Value = 2864.21742
$ws.Range("L138").Value = 16256.3481
$ws.Range("M138").Value = 2275.78258
$ws.Range("N138").Value = -26536.3481

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5414.1084
$ws.Range("I32").Value = 3824.5972
$ws.Range("J32").Value = 15818.182
$ws.Range("K32").Value = 3824.5972
$ws.Range("L32").Value = 15818.182
$ws.Range("M32").Value = -3537.5972
$ws.Range("N32").Value = -16392.182
$ws.Range("H61").Value = 3655.476
$ws.Range("I61").Value = 3766.5789
$ws.Range("J61").Value = 2600
$ws.Range("K61").Value = 3766.5789
$ws.Range("L61").Value = 2600
$ws.Range("M61").Value = -3554.5789
$ws.Range("N61").Value = -3024
$ws.Range("H74").Value = 1425.4333
$ws.Range("I74").Value = 1603.6428
$ws.Range("J74").Value = 1269.5
$ws.Range("K74").Value = 1603.6428
$ws.Range("L74").Value = 1269.5
$ws.Range("M74").Value = -729.6428000000001
$ws.Range("N74").Value = -3017.5
$ws.Range("H77").Value = 1425.4333
$ws.Range("I77").Value = 1603.6428
$ws.Range("J77").Value = 1269.5
$ws.Range("K77").Value = 8018.214
$ws.Range("L77").Value = 6347.5
$ws.Range("M77").Value = -3650.214
$ws.Range("N77").Value = -15083.5
$ws.Range("H102").Value = 2180792
$ws.Range("I102").Value = 2471231.5
$ws.Range("J102").Value = 2495.5
$ws.Range("K102").Value = 2471231.5
$ws.Range("L102").Value = 2495.5
$ws.Range("M102").Value = -2469609.5
$ws.Range("N102").Value = -5739.5
$ws.Range("H132").Value = 2944409.2
$ws.Range("I132").Value = 2647.3809
$ws.Range("J132").Value = 7696486.5
$ws.Range("K132").Value = 7942.1427
$ws.Range("L132").Value = 23089459.5
$ws.Range("M132").Value = -5412.1427
$ws.Range("N132").Value = -23094519.5
$ws.Range("H136").Value = 3655.476
$ws.Range("I136").Value = 3766.5789
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 11299.7367
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -8749.736699999999
$ws.Range("N136").Value = -12900

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 125001560
$ws.Range("I99").Value = 200001040
$ws.Range("J99").Value = 2433.3333
$ws.Range("K99").Value = 200001040
$ws.Range("L99").Value = 2433.3333
$ws.Range("M99").Value = -199999542
$ws.Range("N99").Value = -5429.3333
$ws.Range("H107").Value = 166668420
$ws.Range("I107").Value = 200001900
$ws.Range("K107").Value = 200001900
$ws.Range("M107").Value = -199999980
$ws.Range("H134").Value = 6551.7393
$ws.Range("I134").Value = 9655.846
$ws.Range("J134").Value = 2516.4
$ws.Range("K134").Value = 28967.538
$ws.Range("L134").Value = 7549.200000000001
$ws.Range("M134").Value = -26432.538
$ws.Range("N134").Value = -12619.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 217870.9
$ws.Range("I31").Value = 1662.2609
$ws.Range("J31").Value = 600393.9
$ws.Range("K31").Value = 1662.2609
$ws.Range("L31").Value = 600393.9
$ws.Range("M31").Value = -1367.2609
$ws.Range("N31").Value = -600983.9
$ws.Range("H34").Value = 217870.9
$ws.Range("I34").Value = 1662.2609
$ws.Range("J34").Value = 600393.9
$ws.Range("K34").Value = 1662.2609
$ws.Range("L34").Value = 600393.9
$ws.Range("M34").Value = -1460.2609
$ws.Range("N34").Value = -600797.9
$ws.Range("H58").Value = 1113.7693
$ws.Range("I58").Value = 647.2
$ws.Range("J58").Value = 2669
$ws.Range("K58").Value = 647.2
$ws.Range("L58").Value = 2669
$ws.Range("M58").Value = -444.2
$ws.Range("N58").Value = -3075
$ws.Range("H94").Value = 4117.96
$ws.Range("I94").Value = 5188.1113
$ws.Range("J94").Value = 3516
$ws.Range("K94").Value = 5188.1113
$ws.Range("L94").Value = 3516
$ws.Range("M94").Value = -4737.1113
$ws.Range("N94").Value = -4418
$ws.Range("H99").Value = 9626602
$ws.Range("I99").Value = 12240.5
$ws.Range("J99").Value = 25009580
$ws.Range("K99").Value = 12240.5
$ws.Range("L99").Value = 25009580
$ws.Range("M99").Value = -10742.5
$ws.Range("N99").Value = -25012576
$ws.Range("H107").Value = 20834222
$ws.Range("I107").Value = 33333874
$ws.Range("J107").Value = 1468
$ws.Range("K107").Value = 33333874
$ws.Range("L107").Value = 1468
$ws.Range("M107").Value = -33331954
$ws.Range("N107").Value = -5308
$ws.Range("H126").Value = 9626602
$ws.Range("I126").Value = 12240.5
$ws.Range("J126").Value = 25009580
$ws.Range("K126").Value = 36721.5
$ws.Range("L126").Value = 75028740
$ws.Range("M126").Value = -34251.5
$ws.Range("N126").Value = -75033680
$ws.Range("H132").Value = 1916.06
$ws.Range("I132").Value = 1555.375
$ws.Range("J132").Value = 3358.8
$ws.Range("K132").Value = 4666.125
$ws.Range("L132").Value = 10076.4
$ws.Range("M132").Value = -2136.125
$ws.Range("N132").Value = -15136.4
$ws.Range("H134").Value = 2182.745
$ws.Range("I134").Value = 2640.6128
$ws.Range("J134").Value = 1473.05
$ws.Range("K134").Value = 7921.8384
$ws.Range("L134").Value = 4419.15
$ws.Range("M134").Value = -5386.8384
$ws.Range("N134").Value = -9489.15
$ws.Range("H136").Value = 1113.7693
$ws.Range("I136").Value = 647.2
$ws.Range("J136").Value = 2669
$ws.Range("K136").Value = 1941.6
$ws.Range("L136").Value = 8007
$ws.Range("M136").Value = 608.3999999999999
$ws.Range("N136").Value = -13107
$ws.Range("H141").Value = 33407.555
$ws.Range("J141").Value = 33407.555
$ws.Range("L141").Value = 33407.555
$ws.Range("N141").Value = -43767.555

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 353.5
$ws.Range("J92").Value = 451.25
$ws.Range("L92").Value = 1353.75
$ws.Range("N92").Value = -3849.75
$ws.Range("H122").Value = 5605.5
$ws.Range("J122").Value = 9448.392
$ws.Range("L122").Value = 85035.52799999999
$ws.Range("N122").Value = -89935.52799999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 24775
$ws.Range("J52").Value = 24775
$ws.Range("L52").Value = 24775
$ws.Range("N52").Value = -25293
$ws.Range("H97").Value = 933.3333
$ws.Range("I97").Value = 928.5714
$ws.Range("J97").Value = 950
$ws.Range("K97").Value = 928.5714
$ws.Range("L97").Value = 950
$ws.Range("M97").Value = -432.5714
$ws.Range("N97").Value = -1942
$ws.Range("H102").Value = 1314.4849
$ws.Range("I102").Value = 1012.4167
$ws.Range("J102").Value = 2120
$ws.Range("K102").Value = 1012.4167
$ws.Range("L102").Value = 2120
$ws.Range("M102").Value = 609.5833
$ws.Range("N102").Value = -5364

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12726374
$ws.Range("I132").Value = 17241348
$ws.Range("J132").Value = 2356.2727
$ws.Range("K132").Value = 51724044
$ws.Range("L132").Value = 7068.8181
$ws.Range("M132").Value = -51721514
$ws.Range("N132").Value = -12128.8181
$ws.Range("H136").Value = 7478.952
$ws.Range("I136").Value = 5151.7715
$ws.Range("J136").Value = 19114.857
$ws.Range("K136").Value = 15455.3145
$ws.Range("L136").Value = 57344.571
$ws.Range("M136").Value = -12905.3145
$ws.Range("N136").Value = -62444.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8648.481
$ws.Range("I14").Value = 6818.1816
$ws.Range("J14").Value = 9906.8125
$ws.Range("K14").Value = 6818.1816
$ws.Range("L14").Value = 9906.8125
$ws.Range("M14").Value = -6650.1816
$ws.Range("N14").Value = -10242.8125
$ws.Range("H62").Value = 3375.5
$ws.Range("I62").Value = 3375.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3375.5
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -2751.5
$ws.Range("H65").Value = 3375.5
$ws.Range("I65").Value = 3375.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16877.5
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -13757.5
$ws.Range("H122").Value = 1144.375
$ws.Range("I122").Value = 1139.3182
$ws.Range("K122").Value = 3417.9546
$ws.Range("M122").Value = -967.9546
$ws.Range("H132").Value = 1112.3572
$ws.Range("I132").Value = 806.9655
$ws.Range("J132").Value = 1793.6154
$ws.Range("K132").Value = 2420.8965
$ws.Range("L132").Value = 5380.8462
$ws.Range("M132").Value = 109.1035000000002
$ws.Range("N132").Value = -10440.8462
$ws.Range("H136").Value = 10002872
$ws.Range("I136").Value = 3489.1333
$ws.Range("J136").Value = 25001946
$ws.Range("K136").Value = 10467.3999
$ws.Range("L136").Value = 75005838
$ws.Range("M136").Value = -7917.3999
$ws.Range("N136").Value = -75010938
